$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$ws.Range("BD2").Value = "iahmed@govolution.com"
$ws.Range("BD3").Value = "iahmed@govolution.com"
$ws.Range("BD4").Value = "iahmed@govolution.com"
$ws.Range("BD5").Value = "iahmed@govolution.com"
